# Append the newly-scraped price row (2026-02-07) to the tracking sheet.
# Columns: A=Date, B=Price, C=Discount, D=Incredible.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow   = 38
$dateVal  = "2026-02-07"
$priceVal = "483680"
$discVal  = "41"
$incVal   = "0"

# Every existing cell in this table is stored as plain text (even though the
# values look like dates/numbers). A naive `$ws.Range(...).Value = $dateVal`
# would let Excel auto-detect these as a real date / numbers (like typing
# them into a General cell) and reformat the cells accordingly.
#
# To add them as genuine text - without touching any cell's number format /
# style - compute each value as text via TEXT() in a scratch area well
# outside the used range, then copy only the *values* (not formats) into the
# new row. PasteSpecial(xlPasteValues) brings over the already-evaluated
# text without re-parsing it as user input, so the destination keeps its
# original (default) style while still being stored as text.
$ws.Range("H1").Formula = "=TEXT(""$dateVal"",""@"")"
$ws.Range("H2").Formula = "=TEXT($priceVal,""@"")"
$ws.Range("H3").Formula = "=TEXT($discVal,""@"")"
$ws.Range("H4").Formula = "=TEXT($incVal,""@"")"

$ws.Range("H1").Copy()
$ws.Range("A$newRow").PasteSpecial(-4163)
$ws.Range("H2").Copy()
$ws.Range("B$newRow").PasteSpecial(-4163)
$ws.Range("H3").Copy()
$ws.Range("C$newRow").PasteSpecial(-4163)
$ws.Range("H4").Copy()
$ws.Range("D$newRow").PasteSpecial(-4163)

# Clean up the scratch cells and clipboard marker so nothing else changes.
$ws.Range("H1:H4").ClearContents()
$excel.CutCopyMode = $false
